# Adds a new "2022-Q1" sheet (fund holdings detail) before the "总计"
# (totals) sheet, and inserts a corresponding summary row into "总计".

function Set-TextValue($range, $value) {
    # Force a value to be stored as text (keeps leading zeros / avoids
    # numeric auto-conversion), then drop the temporary "Text" number
    # format so the cell ends up with no special style, same as a
    # plain inline string cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# "总计" is currently the last sheet (index 4).
$totalIndex = $wb.Worksheets.Count
$totalSheet = $wb.Worksheets.Item($totalIndex)

# Duplicate "总计" to inherit its formatting/styles, placing the copy
# right before it; this copy becomes the new "2022-Q1" sheet. Note:
# after Copy(), the new duplicate ends up at $totalIndex and the
# original "总计" is pushed to $totalIndex + 1, so re-fetch both by
# position instead of trusting old object references.
$totalSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item($totalIndex)
$totalSheet = $wb.Worksheets.Item($totalIndex + 1)
$newSheet.Name = "2022-Q1"

# ---- Header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
# Give the newly added E1:H1 header cells the same style as B1:D1.
$newSheet.Range("B1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

# ---- Data rows ----
Set-TextValue $newSheet.Range("B2") "210009"
Set-TextValue $newSheet.Range("C2") "金鹰核心资源混合"
Set-TextValue $newSheet.Range("D2") "3.86"
Set-TextValue $newSheet.Range("E2") "94.96"
Set-TextValue $newSheet.Range("F2") "4.74"
Set-TextValue $newSheet.Range("G2") "0.1830"
$newSheet.Range("H2").Value = 6

Set-TextValue $newSheet.Range("B3") "162102"
Set-TextValue $newSheet.Range("C3") "金鹰中小盘精选混合"
Set-TextValue $newSheet.Range("D3") "4.60"
Set-TextValue $newSheet.Range("E3") "76.52"
Set-TextValue $newSheet.Range("F3") "3.78"
Set-TextValue $newSheet.Range("G3") "0.1739"
$newSheet.Range("H3").Value = 4

Set-TextValue $newSheet.Range("B4") "001167"
Set-TextValue $newSheet.Range("C4") "金鹰科技创新股票"
Set-TextValue $newSheet.Range("D4") "4.03"
Set-TextValue $newSheet.Range("E4") "94.55"
Set-TextValue $newSheet.Range("F4") "4.31"
Set-TextValue $newSheet.Range("G4") "0.1737"
$newSheet.Range("H4").Value = 10

# ---- Update "总计" with the new 2022-Q1 summary row ----
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.53

# Restore A2's style (row-insert formatting changed it); A4 still has
# the original style to copy from.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

# The shifted rows' index column (A) values need to be bumped by one
# to keep representing a sequential 0-based row index.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
